# Update attendance / price figures on the "展览" and "全部类型" sheets.
# Both sheets contain the same table of convention listings, so the same
# cell updates are applied to each.

$wb = $excel.ActiveWorkbook

$targetSheetNames = @("展览", "全部类型")

# Map of cell address -> new value for each affected sheet.
$updates = @{
    "F5"  = 2561
    "F9"  = 1308
    "F12" = 13
    "F13" = 1157
    "G13" = 55
    "F14" = 335
    "F15" = 322
    "F16" = 28
    "F21" = 2345
    "F23" = 271
}

foreach ($sheetName in $targetSheetNames) {
    $ws = $wb.Worksheets.Item($sheetName)
    foreach ($addr in $updates.Keys) {
        $ws.Range($addr).Value = $updates[$addr]
    }
}
